$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the two observation records currently stored in rows 13
# (LC / Dropptaggsvamp) and 14 (NT / Spillkråka): row 13 becomes the old
# row 14 record and vice versa. Only the columns that actually differ
# between the two records are touched; columns that already share the
# same value in both rows (P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AY)
# are left untouched.

# --- Capture current (pre-edit) values for the columns that differ ---
$A13 = $ws.Range("A13").Value2()
$A14 = $ws.Range("A14").Value2()

$B13 = $ws.Range("B13").Value2()
$B14 = $ws.Range("B14").Value2()

$D13 = $ws.Range("D13").Value2()
$D14 = $ws.Range("D14").Value2()

$E13 = $ws.Range("E13").Value2()
$E14 = $ws.Range("E14").Value2()

$F13 = $ws.Range("F13").Value2()
$F14 = $ws.Range("F14").Value2()

$G13 = $ws.Range("G13").Value2()
$G14 = $ws.Range("G14").Value2()

$H13 = $ws.Range("H13").Value2()
$H14 = $ws.Range("H14").Value2()

$I13 = $ws.Range("I13").Value2()
$I14 = $ws.Range("I14").Value2()

$J13 = $ws.Range("J13").Value2()
$J14 = $ws.Range("J14").Value2()

$K13 = $ws.Range("K13").Value2()
$K14 = $ws.Range("K14").Value2()

$M13 = $ws.Range("M13").Value2()
$M14 = $ws.Range("M14").Value2()

$Q13 = $ws.Range("Q13").Value2()
$Q14 = $ws.Range("Q14").Value2()

$R13 = $ws.Range("R13").Value2()
$R14 = $ws.Range("R14").Value2()

$Z13 = $ws.Range("Z13").Value2()
$Z14 = $ws.Range("Z14").Value2()

$AB13 = $ws.Range("AB13").Value2()
$AB14 = $ws.Range("AB14").Value2()

$AX13 = $ws.Range("AX13").Value2()
$AX14 = $ws.Range("AX14").Value2()

# --- Numeric columns: plain numbers, swap directly ---
$ws.Range("A13").Value = $A14
$ws.Range("A14").Value = $A13

$ws.Range("B13").Value = $B14
$ws.Range("B14").Value = $B13

$ws.Range("E13").Value = $E14
$ws.Range("E14").Value = $E13

$ws.Range("Q13").Value = $Q14
$ws.Range("Q14").Value = $Q13

$ws.Range("R13").Value = $R14
$ws.Range("R14").Value = $R13

# --- Plain text columns: swap directly ---
$ws.Range("D13").Value = $D14
$ws.Range("D14").Value = $D13

$ws.Range("F13").Value = $F14
$ws.Range("F14").Value = $F13

$ws.Range("G13").Value = $G14
$ws.Range("G14").Value = $G13

$ws.Range("H13").Value = $H14
$ws.Range("H14").Value = $H13

$ws.Range("Z13").Value = $Z14
$ws.Range("Z14").Value = $Z13

$ws.Range("AB13").Value = $AB14
$ws.Range("AB14").Value = $AB13

$ws.Range("AX13").Value = $AX14
$ws.Range("AX14").Value = $AX13

# --- J / K / M: present on only one of the two rows, so clear first
#     then (re)populate as needed instead of a pure swap ---
$ws.Range("J13").ClearContents()
$ws.Range("J14").ClearContents()
$ws.Range("K13").ClearContents()
$ws.Range("K14").ClearContents()
$ws.Range("M13").ClearContents()
$ws.Range("M14").ClearContents()

if ($J13 -ne $null) { $ws.Range("J14").Value = $J13 }
if ($J14 -ne $null) { $ws.Range("J13").Value = $J14 }

if ($K13 -ne $null) { $ws.Range("K14").Value = $K13 }
if ($K14 -ne $null) { $ws.Range("K13").Value = $K14 }

if ($M13 -ne $null) { $ws.Range("M14").Value = $M13 }
if ($M14 -ne $null) { $ws.Range("M13").Value = $M14 }

# --- Column I holds text that looks like a plain integer ("5", "1").
#     Prefix with an apostrophe (same as typing '5 into Excel) so the
#     stored type stays Text (matching the original inlineStr cells)
#     instead of being inferred as a number, without touching NumberFormat. ---
$ws.Range("I13").Value = "'" + $I14
$ws.Range("I14").Value = "'" + $I13
